# "fixed error in experiment format"
#
# The "Survey 2" sheet was missing the "Pseudo-Random Question Width:"
# row that "Survey 1" already has in its [Survey Table Properties]
# block (row 8, right under "Questions Per Page:"). Insert the missing
# row so both survey sheets share the same layout, then restore the
# normal (un-navigated) view state: Survey 2 becomes the active sheet,
# with the newly inserted row selected, and Survey 1 loses the stale
# "currently selected tab" / scroll-position markers left over from
# editing.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Survey 2")

# Insert the missing "Pseudo-Random Question Width:" row into Survey 2,
# shifting rows 8:68 down to 9:69 (mirrors Survey 1's existing row 8).
$ws2.Rows.Item(8).Insert() | Out-Null
$ws2.Range("A8").Value = "Pseudo-Random Question Width:"

# Restore view state: Survey 2 active, new row selected; Survey 1's
# leftover tab-selection/scroll position cleared (selection untouched).
$ws2.Activate() | Out-Null
$ws2.Range("A8:D8").Select() | Out-Null
